# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) for specific rows across all 8 profession sheets, per the
# scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 410.45456
$ws.Range("J17").Value = 451.7143
$ws.Range("L17").Value = 1355.1429
$ws.Range("N17").Value = -1691.1429
$ws.Range("H76").Value = 7421.067
$ws.Range("I76").Value = 4268
$ws.Range("K76").Value = 4268
$ws.Range("M76").Value = -3953
$ws.Range("H79").Value = 7421.067
$ws.Range("I79").Value = 4268
$ws.Range("K79").Value = 4268
$ws.Range("M79").Value = -3176
$ws.Range("H88").Value = 915396.6
$ws.Range("J88").Value = 1117272.5
$ws.Range("L88").Value = 1117272.5
$ws.Range("N88").Value = -1118084.5
$ws.Range("H91").Value = 915396.6
$ws.Range("J91").Value = 1117272.5
$ws.Range("L91").Value = 1117272.5
$ws.Range("N91").Value = -1120080.5
$ws.Range("H113").Value = 9514.556
$ws.Range("I113").Value = 13881
$ws.Range("K113").Value = 13881
$ws.Range("M113").Value = -10627
$ws.Range("H118").Value = 598.44446
$ws.Range("I118").Value = 598.44446
$ws.Range("K118").Value = 1795.33338
$ws.Range("M118").Value = -138.33338
$ws.Range("H137").Value = 2449.7222
$ws.Range("I137").Value = 2643.6428
$ws.Range("K137").Value = 7930.928400000001
$ws.Range("M137").Value = -5380.928400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 365.5263
$ws.Range("I97").Value = 314.57144
$ws.Range("J97").Value = 508.2
$ws.Range("K97").Value = 314.57144
$ws.Range("L97").Value = 508.2
$ws.Range("M97").Value = 181.42856
$ws.Range("N97").Value = -1500.2
$ws.Range("H132").Value = 7383.65
$ws.Range("I132").Value = 3801.2903
$ws.Range("J132").Value = 19722.889
$ws.Range("K132").Value = 11403.8709
$ws.Range("L132").Value = 59168.667
$ws.Range("M132").Value = -8873.8709
$ws.Range("N132").Value = -64228.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 731.6667
$ws.Range("I16").Value = 731.6667
$ws.Range("K16").Value = 731.6667
$ws.Range("M16").Value = -561.6667
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H132").Value = 98493.5
$ws.Range("J132").Value = 98493.5
$ws.Range("L132").Value = 98493.5
$ws.Range("N132").Value = -108613.5
$ws.Range("H135").Value = 78916.586
$ws.Range("J135").Value = 78916.586
$ws.Range("L135").Value = 78916.586
$ws.Range("N135").Value = -89056.586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1233.04
$ws.Range("I31").Value = 1306.5454
$ws.Range("J31").Value = 1175.2858
$ws.Range("K31").Value = 1306.5454
$ws.Range("L31").Value = 1175.2858
$ws.Range("M31").Value = -1011.5454
$ws.Range("N31").Value = -1765.2858
$ws.Range("H34").Value = 1233.04
$ws.Range("I34").Value = 1306.5454
$ws.Range("J34").Value = 1175.2858
$ws.Range("K34").Value = 1306.5454
$ws.Range("L34").Value = 1175.2858
$ws.Range("M34").Value = -1104.5454
$ws.Range("N34").Value = -1579.2858
$ws.Range("H58").Value = 2316.3125
$ws.Range("I58").Value = 2287.7273
$ws.Range("J58").Value = 2379.2
$ws.Range("K58").Value = 2287.7273
$ws.Range("L58").Value = 2379.2
$ws.Range("M58").Value = -2084.7273
$ws.Range("N58").Value = -2785.2
$ws.Range("H99").Value = 3246.147
$ws.Range("I99").Value = 3050.3462
$ws.Range("J99").Value = 3882.5
$ws.Range("K99").Value = 3050.3462
$ws.Range("L99").Value = 3882.5
$ws.Range("M99").Value = -1552.3462
$ws.Range("N99").Value = -6878.5
$ws.Range("H107").Value = 1350.25
$ws.Range("I107").Value = 1001.2
$ws.Range("K107").Value = 1001.2
$ws.Range("M107").Value = 918.8
$ws.Range("H126").Value = 3246.147
$ws.Range("I126").Value = 3050.3462
$ws.Range("J126").Value = 3882.5
$ws.Range("K126").Value = 9151.0386
$ws.Range("L126").Value = 11647.5
$ws.Range("M126").Value = -6681.0386
$ws.Range("N126").Value = -16587.5
$ws.Range("H136").Value = 2316.3125
$ws.Range("I136").Value = 2287.7273
$ws.Range("J136").Value = 2379.2
$ws.Range("K136").Value = 6863.1819
$ws.Range("L136").Value = 7137.599999999999
$ws.Range("M136").Value = -4313.1819
$ws.Range("N136").Value = -12237.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 178.22223
$ws.Range("I40").Value = 120.15385
$ws.Range("K40").Value = 480.6154
$ws.Range("M40").Value = -411.6154
$ws.Range("H59").Value = 1999.6666
$ws.Range("I59").Value = 1999
$ws.Range("J59").Value = 1999.8
$ws.Range("K59").Value = 5997
$ws.Range("L59").Value = 5999.4
$ws.Range("M59").Value = -5457
$ws.Range("N59").Value = -7079.4
$ws.Range("I68").Value = 301
$ws.Range("J68").Value = 999
$ws.Range("K68").Value = 903
$ws.Range("L68").Value = 2997
$ws.Range("M68").Value = -92
$ws.Range("N68").Value = -4619
$ws.Range("I71").Value = 301
$ws.Range("J71").Value = 999
$ws.Range("K71").Value = 2709
$ws.Range("L71").Value = 8991
$ws.Range("M71").Value = 1347
$ws.Range("N71").Value = -17103

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9324.777
$ws.Range("I5").Value = 2250.75
$ws.Range("J5").Value = 14984
$ws.Range("K5").Value = 2250.75
$ws.Range("L5").Value = 14984
$ws.Range("M5").Value = -2138.75
$ws.Range("N5").Value = -15208
$ws.Range("H97").Value = 655.6
$ws.Range("I97").Value = 449.66666
$ws.Range("J97").Value = 964.5
$ws.Range("K97").Value = 449.66666
$ws.Range("L97").Value = 964.5
$ws.Range("M97").Value = 46.33334000000002
$ws.Range("N97").Value = -1956.5
$ws.Range("H132").Value = 3084.8333
$ws.Range("I132").Value = 2701.8
$ws.Range("K132").Value = 8105.400000000001
$ws.Range("M132").Value = -5575.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7854.2
$ws.Range("J7").Value = 8934.416999999999
$ws.Range("L7").Value = 8934.416999999999
$ws.Range("N7").Value = -9158.416999999999
$ws.Range("H43").Value = 18000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 18000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 18000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -18386
$ws.Range("H93").Value = 1109.7587
$ws.Range("I93").Value = 915.381
$ws.Range("K93").Value = 915.381
$ws.Range("M93").Value = 332.619
$ws.Range("H126").Value = 7854.2
$ws.Range("J126").Value = 8934.416999999999
$ws.Range("L126").Value = 26803.251
$ws.Range("N126").Value = -31743.251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5998.6665
$ws.Range("J2").Value = 4996
$ws.Range("L2").Value = 4996
$ws.Range("N2").Value = -5220
$ws.Range("H126").Value = 3437.1538
$ws.Range("I126").Value = 3531.9167
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 10595.7501
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -8125.750100000001
$ws.Range("N126").Value = -11840
